$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.801.69"
$ws.Range("E2").Value = "  +2.91%  "

$ws.Range("D3").Value = "1.912.57"
$ws.Range("E3").Value = "  +2.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.029"
$ws.Range("E4").Value = "  +2.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.02"
$ws.Range("E5").Value = "  +3.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.033"
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5222"
$ws.Range("E7").Value = "  +1.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3972"
$ws.Range("E8").Value = "  +3.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08391"
$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.139"
$ws.Range("E10").Value = "  +2.71%  "

$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.922.45"
$ws.Range("E11").Value = "  +3.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.307"
$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.73"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.355"
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.026"
$ws.Range("E15").Value = "  +2.26%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001118"
$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.78"
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06809"
$ws.Range("E18").Value = "  +2.34%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.10"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.028"
$ws.Range("E20").Value = "  +2.54%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.112"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "28.823.53"
$ws.Range("E22").Value = "  +2.87%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.33"
$ws.Range("E23").Value = "  +2.41%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.298"
$ws.Range("E24").Value = "  +2.29%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.435"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("D26").Value = "2.140.14"
$ws.Range("E26").Value = "  +3.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.71"
$ws.Range("E27").Value = "  +3.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.02"
$ws.Range("E28").Value = "  +2.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.478"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.83"
$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1062"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.053"
$ws.Range("E32").Value = "  +2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.968"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.693"
$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.538"
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02494"
$ws.Range("E36").Value = "  +3.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06662"
$ws.Range("E37").Value = "  +2.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2243"
$ws.Range("E38").Value = "  +3.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6596"
$ws.Range("E39").Value = "  +0.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.265"
$ws.Range("E40").Value = "  +3.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.201"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.052"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.13"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6192"
$ws.Range("E44").Value = "  +0.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.18"
$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.780"
$ws.Range("E46").Value = "  +3.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.309"
$ws.Range("E47").Value = "  +2.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.035"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.246"
$ws.Range("E49").Value = "  +2.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.20"
$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06957"
$ws.Range("E51").Value = "  +1.71%  "
